# Generate Report for Handback
#
# The handback step for the "91046edf-...md" file failed because the
# generated handback file name didn't match the original handoff file
# name. Reflect that on both locale report sheets:
#   - Status (column C, row 7) flips from "Ready for handoff" to
#     "Handback transform failed" (this text lives in the shared string
#     table, so the Overview roll-up sheet picks the change up too).
#   - Error Detail (column P, row 7) gets a locale-specific explanation
#     of the filename mismatch.
#   - The Error Detail column is widened so the longer message is legible.

$wb = $excel.ActiveWorkbook

# Update the status text everywhere it appears (Overview, zh-cn, de-de) so
# the shared underlying text changes in place rather than leaving the old
# wording behind on sheets this script doesn't touch directly.
foreach ($sheet in $wb.Worksheets) {
    [void]$sheet.Cells.Replace("Ready for handoff", "Handback transform failed")
}

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Columns.Item(16).ColumnWidth = 39.17
$wsZh.Range("P7").Value = "Handback file name: oxckfawo.fds is different with handoff file name: 91046edf-a6f8-40cf-8cfc-d35936c20c05.20f231db9a543960a61df7d960b625b0335792a8.zh-cn."

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Columns.Item(16).ColumnWidth = 39.17
$wsDe.Range("P7").Value = "Handback file name: oxckfawo.fds is different with handoff file name: 91046edf-a6f8-40cf-8cfc-d35936c20c05.20f231db9a543960a61df7d960b625b0335792a8.de-de."
